$wb = $excel.ActiveWorkbook

# --- Sheet: "døgnrytmetabell" (rId5) -- update the activity-count grid (rows 26-49, cols C:I) ---
$ws = $wb.Worksheets.Item("døgnrytmetabell")

$ws.Range("C26").Value = 3
$ws.Range("C28:I28").Value = 3
$ws.Range("C29:I29").Value = 3
$ws.Range("C30:I30").Value = 3
$ws.Range("C31:I31").Value = 4
$ws.Range("C32:I32").Value = 5
$ws.Range("C33:I33").Value = 6
$ws.Range("C34:I34").Value = 7
$ws.Range("C35:I35").Value = 7
$ws.Range("C36:I36").Value = 8
$ws.Range("C37:I37").Value = 8
$ws.Range("H38:I38").Value = 9
$ws.Range("G41").Value = 9
$ws.Range("G42").Value = 9
$ws.Range("C43:I43").Value = 9
$ws.Range("G43").Value = 8
$ws.Range("C44:I44").Value = 9
$ws.Range("G44").Value = 7
$ws.Range("G45").Value = 6
$ws.Range("C46:I46").Value = 7
$ws.Range("G46").Value = 6
$ws.Range("C47:I47").Value = 6
$ws.Range("C48:I48").Value = 5
$ws.Range("I49").Value = 4

# --- Update sheet selections / active-tab state to match the saved UI state ---

# "døgnrytmetabell (2)" -- selection moved to A13 (single cell)
$wsA = $wb.Worksheets.Item("døgnrytmetabell (2)")
$wsA.Range("A13").Select()

# "bemanningsplan" -- stays selected at M14, but loses tab focus (scrolled, no longer the active tab)
$wsB = $wb.Worksheets.Item("bemanningsplan")
$wsB.Range("M14").Select()

# "døgnrytmetabell" -- ends up the active/visible tab, selection moved to F52
$ws.Range("F52").Select()
